$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append transaction rows 299-381 (new batch synced 2026-02-10)
$data = @(
    ,@("2026-02-10 12:07:35", "237681102046", "CLARION MENKE BAKARI", 79120)
    ,@("2026-02-10 14:49:18", "237681114247", "MINDEM SARL ymele voufoYMELE VOUFO VANESSA", 1014)
    ,@("2026-02-10 14:18:14", "237681114370", "ETS LE CONTENT SEDJINE TCHINDA RAMSES GAUTIER", 1847)
    ,@("2026-02-10 12:23:50", "237681118330", "SAHA NDESA JONAS LTDLA_POLAS_OTH_NDOGBONG SERIE", 408039)
    ,@("2026-02-10 17:18:23", "237681125655", "EMENGUE PICHOU ROMEO KAMILAH CONNECTION GROUP", 137219)
    ,@("2026-02-10 16:31:19", "237681180496", "ALEX NGOUO YOUNDA", 9944)
    ,@("2026-02-10 16:00:49", "237681240793", "MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS", 2947)
    ,@("2026-02-10 14:43:17", "237681299829", "NDEBI MEDARD DESIRE ETS MOBILE FINANCIAL SERVICES MFS", 698)
    ,@("2026-02-10 13:18:12", "237681446273", "CHRISTIANE MARTINE ALEXANDRINE NSANG EPSE ESSAKA EBOUMBOU", 43120)
    ,@("2026-02-10 16:48:56", "237681490029", "MOISE LONTCHI", 7460)
    ,@("2026-02-10 12:14:45", "237681589841", "LEINTENG ROSE MARY", 70456)
    ,@("2026-02-10 15:07:51", "237681602244", "TSOMEJIO KENFACK NICAISE NESLIE ETS TCHATCHOUANG PAUL  ETP", 62339)
    ,@("2026-02-10 12:34:52", "237681603496", "ADVINE STEPHANIE NGOUNGO WABEU", 295896)
    ,@("2026-02-10 16:01:21", "237681606646", "DERRICK SONWA LONTIO", 114713)
    ,@("2026-02-10 15:26:57", "237681611433", "Barry Diakariaou World T Plus", 79704)
    ,@("2026-02-10 16:04:43", "237681655237", "ETS LE CONTENT DJOUFACK WOUAFACK ALAIN GAROUSTE", 46613)
    ,@("2026-02-10 15:42:15", "237681655241", "LA NEGRESSE LTDLA CBOX R1 NKUIDJEU KAMDOUM SYMPHORIEN", 18433)
    ,@("2026-02-10 17:53:01", "237681656314", "SWIRRI AZINWI NGANG LA NEGRESSE SARL", 481905)
    ,@("2026-02-10 16:42:32", "237681657461", "LA NEGRESSE SARL MOKAM KOUAM VIVIANE", 215288)
    ,@("2026-02-10 15:41:00", "237681657562", "BLANDINE PEYEMBOUO", 81593)
    ,@("2026-02-10 15:39:44", "237681657939", "ETS MOBILE FINANCIAL SERVICES MFS LTDLA CBOX R1 MOHA CHAIBOU", 1)
    ,@("2026-02-10 16:25:01", "237681658403", "LA NEGRESSE MISSOKE-UNIVERSITE", 31000)
    ,@("2026-02-10 13:09:33", "237681662596", "LA NEGRESSE SARL LTDLA CBOX R1 TIOKENG SANDRINE", 589452)
    ,@("2026-02-10 16:45:41", "237681662606", "GAEL PHALENNE NANA POUASSI", 119513)
    ,@("2026-02-10 15:05:40", "237681662680", "pascaline djiogo mejioguezem", 112445)
    ,@("2026-02-10 15:44:12", "237681662761", "ETS AMOUR DE DIEU SERVICES LTDLA_CBOX_R1_MBOCK NICOLE RUTH", 61980)
    ,@("2026-02-10 14:41:55", "237681663743", "LA NEGRESSE SARL FONGA SINTCHA YOLANDE MIREILLE", 441531)
    ,@("2026-02-10 16:33:52", "237681676445", "MELANIE NGAFFO", 401091)
    ,@("2026-02-09 15:12:33", "237681677617", "ETS MOBILE FINANCIAL SERVICES MFS MANFOUO TCHOUALA HUGUETTE", 79864)
    ,@("2026-02-10 15:30:17", "237681678622", "FOKGO BRIGITTE ETS MOBILE FINANCIAL SERVICES MFS", 122195)
    ,@("2026-02-10 21:51:16", "237681679096", "TAMNOU NGANGO ULRICH BERNARD ETS MOBILE FINANCIAL SERVICES MFS", 124569)
    ,@("2026-02-10 17:06:49", "237681679214", "TIENTCHEU ROSINE CHRISTELLE ETS MOBILE FINANCIAL SERVICES MFS", 26962)
    ,@("2026-02-10 18:14:48", "237681679310", "LA NEGRESSE SARL DIALLO AMADOU OURY", 1929)
    ,@("2026-02-10 14:16:02", "237681679880", "TSAKEM AGNES LAIDY ETS LE CONTENT", 18717)
    ,@("2026-02-10 15:02:38", "237681862876", "TIDO GARLINE NOGRA-POLAS-BTQ-MAKEPE MISSOKE", 10280)
    ,@("2026-02-10 16:07:12", "237682117915", "MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS", 132271)
    ,@("2026-02-10 13:18:35", "237682154553", "N A ISUFUH MIEMONA NGESSY ETS MOBILE FINANCIAL SERVICES MFS", 20837)
    ,@("2026-02-10 14:41:36", "237682238745", "RACHEL PRUDENCE JIKE KETCHA", 734434)
    ,@("2026-02-10 12:41:01", "237682316602", "DANGA ZAMPA PATRICE VICKY TOP MOBIL TELECOM", 898713)
    ,@("2026-02-10 09:53:22", "237682323406", "ETS LE CONTENT LAS VEGAS", 26)
    ,@("2026-02-10 16:59:09", "237682368679", "MFS SIM AA 2", 24203)
    ,@("2026-02-10 14:20:36", "237682370358", "CARINE SONKENG", 350161)
    ,@("2026-02-10 14:15:47", "237682430965", "ELSA CABRELLE MAKOUNGANG ETS MOBILE FINANCIAL SERVICES MFS", 57098)
    ,@("2026-02-10 13:20:24", "237682480811", "KENGNE TADJO LYNDA NOEL ETS MOBILE FINANCIAL SERVICES MFS", 597064)
    ,@("2026-02-10 12:38:21", "237682511457", "FRANFORETTE NWOGUEP KODJOUO", 55164)
    ,@("2026-02-10 14:28:27", "237682520113", "KEUYAP NGATCHEU JUDITH JOSY ETS MOBILE FINANCIAL SERVICES MFS", 254361)
    ,@("2026-02-10 15:43:10", "237682639044", "JOSEPH KAMGA", 53913)
    ,@("2026-02-10 14:54:44", "237682764368", "LA NEGRESSE SARL LIEDJI GINETTE", 503106)
    ,@("2026-02-10 14:30:43", "237682798275", "NGAFFO YOCADINE BENEDITE ETS MOBILE FINANCIAL SERVICES MFS", 8147)
    ,@("2026-02-10 02:27:05", "237682803277", "NGUEMASSOM RENE MARTIAL LA NEGRESSE SARL", 17)
    ,@("2026-02-10 14:36:45", "237682814055", "DIALL BOUCARI KAMILAH CONNECTION GROUP", 23560)
    ,@("2026-02-10 15:47:09", "237682827350", "ALAIN GACIEN DOUANLA", 441750)
    ,@("2026-02-10 15:53:40", "237682975726", "LA NEGRESSE SARL NYOUNG JOSEPH CLOTAIRE", 182156)
    ,@("2026-02-09 14:50:06", "237683023087", "FAGHUIE ABIBA", 580555)
    ,@("2026-02-10 17:34:21", "237683075075", "PARFAIT TEMOH DAH", 322352)
    ,@("2026-02-10 17:17:36", "237683079541", "THIERRY MELINGUI AYINA", 200040)
    ,@("2026-02-10 17:18:08", "237683165199", "DIALLO MAMADOU LAMINE ETS MOBILE FINANCIAL SERVICES MFS", 1756)
    ,@("2026-02-10 14:09:38", "237683279255", "Deuffi Yvonne laurette LENA GLOBAL", 3269)
    ,@("2026-02-10 15:11:32", "237683323481", "ETS LE CONTENT TSAZE DONFOUET FLORETTE ROSINE", 576123)
    ,@("2026-02-10 21:05:50", "237683353137", "VOUGMO NGUEMO MERLIN WILLIAM ETS MOBILE FINANCIAL SERVICES MFS", 238659)
    ,@("2026-02-10 15:47:31", "237683356603", "MFS MATCHINDA SENDRINE", 1099)
    ,@("2026-02-10 13:38:08", "237683356768", "RUSSEL LECLER KOUTJEM", 219074)
    ,@("2026-02-10 00:56:11", "237683360459", "LUCIE MAJOLIE LELE NKANKEU", 326)
    ,@("2026-02-10 14:58:43", "237683366333", "ETS MOBILE FINANCIAL SERVICES MFS DJIAGUE JEANNETTE", 33240)
    ,@("2026-02-10 15:54:52", "237683368985", "MFS BELL HENRIE BERNARD", 240076)
    ,@("2026-02-10 15:09:57", "237683379070", "MELI DOUANLA ORNELA LINDA ETS MOBILE FINANCIAL SERVICES MFS", 1309653)
    ,@("2026-02-10 13:03:30", "237683379155", "NIMBUNG EPSE BWEH ODETTE", 739806)
    ,@("2026-02-10 13:23:20", "237683379207", "NANFACK EPSE SOKENG SOLANGE ETS MOBILE FINANCIAL SERVICES MFS", 162662)
    ,@("2026-02-10 16:11:43", "237683386020", "LA NEGRESSE SARL LTDLA CBOX R1 MAFFO DALLY DIANE", 8632)
    ,@("2026-02-10 16:07:47", "237683394976", "DIALLO MAMADOU OURY", 15860)
    ,@("2026-02-10 15:29:50", "237683395123", "LA NEGRESSE SARL LTDLA-CBOX-R1-TAGNIN NICAISSE FLEURIE", 87704)
    ,@("2026-02-10 14:48:08", "237683396173", "ENOMA NDJAH PAULINE SONIA ETS MOBILE FINANCIAL SERVICES MFS", 115623)
    ,@("2026-02-10 13:42:30", "237683400719", "CHANCELINE LAGMAGO", 169858)
    ,@("2026-02-10 13:27:10", "237683408221", "KOMI GISELE ETS MOBILE FINANCIAL SERVICES MFS", 69841)
    ,@("2026-02-10 16:08:28", "237683432110", "ARNAUD GHISLAIN FOSSO", 818170)
    ,@("2026-02-10 13:43:57", "237683454059", "ELISABETH MARIE ETIENNE ANZIN", 59951)
    ,@("2026-02-10 12:35:10", "237683454060", "RTS BP CITÉ", 6635)
    ,@("2026-02-10 14:02:10", "237683555873", "CHI MERCY SWIRI LTDLA_POLAS_BTQ_LIMBE", 1291825)
    ,@("2026-02-10 14:38:12", "237683557193", "KOUYEKE MONIQUE LA NEGRESSE SARL", 22845)
    ,@("2026-02-10 10:14:32", "237683612202", "ALBERTINE TIBELLE DONGMO NANFACK", 15004)
    ,@("2026-02-10 15:31:46", "237683730580", "LA NEGRESSE SARL LTDLA_CBOX_R1_TSAFO NICOLE", 330575)
    ,@("2026-02-10 15:46:49", "237683743490", "ETS LE CONTENT NGAH MARIE", 10613)
    ,@("2026-02-10 14:50:52", "237683815311", "ETS TIN-GLOBALCOMM ZEBAZE TSEBAZE LAURA", 2759)
)

$startRow = 299
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = "'" + $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
